$wb = $excel.ActiveWorkbook

# New field-data rows collected today (2020-10-30) for each site sheet.
# Columns: A = datetime (serial), B = chloride_mgL, C = cl_temp_C
$newRows = @(
    @{ Sheet = "WIC";  Row = 8;  DateTime = 44134.444444444445; Chloride = 90;                 Temp = 4.5 },
    @{ Sheet = "YS";   Row = 21; DateTime = 44134.479166666664; Chloride = 37;                 Temp = 8.4 },
    @{ Sheet = "SW";   Row = 21; DateTime = 44134.525694444441; Chloride = 108;                Temp = 6.5 },
    @{ Sheet = "YI";   Row = 21; DateTime = 44134.325694444444; Chloride = 38.450000000000003; Temp = 7.6 },
    @{ Sheet = "YN";   Row = 21; DateTime = 44134.345138888886; Chloride = 27.74;              Temp = 3.4 },
    @{ Sheet = "6MC";  Row = 21; DateTime = 44134.35833333333;  Chloride = 51.45;              Temp = 4.2 },
    @{ Sheet = "DC";   Row = 21; DateTime = 44134.368055555555; Chloride = 51.6;               Temp = 4.8 },
    @{ Sheet = "PBMS"; Row = 21; DateTime = 44134.384027777778; Chloride = 134;                Temp = 5.2 },
    @{ Sheet = "PBSF"; Row = 21; DateTime = 44134.397222222222; Chloride = 288.2;              Temp = 5.4 }
)

foreach ($entry in $newRows) {
    $ws = $wb.Worksheets.Item($entry.Sheet)
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.DateTime
    $ws.Cells.Item($r, 2).Value = $entry.Chloride
    $ws.Cells.Item($r, 3).Value = $entry.Temp
    $ws.Range("A$r").Select() | Out-Null
}

# PBSF is where data entry finished today, so it ends up the active/visible tab.
$wb.Worksheets.Item("PBSF").Activate() | Out-Null
$wb.Worksheets.Item("PBSF").Range("A21").Select() | Out-Null
